# Update attendance-interest counts ("想去人数") on the "展览" and "全部类型"
# sheets to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F8").Value = 2253   # was 2250
$wsExhibition.Range("F10").Value = 5554  # was 5551
$wsExhibition.Range("F11").Value = 129   # was 128

# --- Sheet "全部类型" ---
$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F11").Value = 2253    # was 2250
$wsAllTypes.Range("F13").Value = 5554    # was 5551
$wsAllTypes.Range("F14").Value = 129     # was 128

$wb.Save()
